$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the paragraph "It will ask for the input file, and type in
# "test.txt"" entirely (it sits right after the first screenshot image).
# ---------------------------------------------------------------------------
$q1 = [char]0x201C
$q2 = [char]0x201D

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "It will ask for the input file*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# ---------------------------------------------------------------------------
# Change 2: extend "The graph is slightly messy, but you can interact with
# it." with a new sentence, then add several new paragraphs describing the
# screenshots / bar graph that follow. The very last paragraph in the body
# is the (empty) one holding the "_GoBack" bookmark; Word recreates it on
# every save, and it stays last throughout these edits.
# ---------------------------------------------------------------------------
$messyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "The graph is slightly messy*") {
        $messyPara = $p
        break
    }
}

# 2a. Append " Use the scroll bars to move around." to that same paragraph.
$messyPara.Range.InsertAfter(" Use the scroll bars to move around.")

# 2b. Insert a new paragraph right after it:
#     There is an image of this graph under drawings folder called
#     "interactive_graph" (with spell-check markers around the bare word,
#     matching how Word tags camel/underscore identifiers it doesn't know).
$insertPoint = $d.Range($messyPara.Range.End - 1, $messyPara.Range.End - 1)
$xml1 = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r><w:t>There is an image of this graph under drawings folder called ' + $q1 + '</w:t></w:r>' +
              '<w:proofErr w:type="spellStart"/>' +
              '<w:r><w:t>interactive_graph</w:t></w:r>' +
              '<w:proofErr w:type="spellEnd"/>' +
              '<w:r><w:t>' + $q2 + '</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$insertPoint.InsertXML($xml1)

# 2c. The last (empty) paragraph of the body holds the "_GoBack" bookmark.
#     Turn it into: "A bar graph ... page ranks" <bookmark> ". You may have
#     to enlarge the bar graph to see the names." The bookmark is a
#     zero-width marker, so re-resolve its position after each edit (rather
#     than trust paragraph/range Start/End, which snap past it) and always
#     insert immediately before it with InsertBefore.
$n = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($n)

$bm = $d.Bookmarks.Item("_GoBack")
$before = $d.Range($bm.Start, $bm.Start)
$before.InsertBefore("A bar graph will show up and present the classes and their perspective page ranks")

$bm2 = $d.Bookmarks.Item("_GoBack")
$after = $d.Range($bm2.Start, $bm2.Start)
$after.InsertBefore(". You may have to enlarge the bar graph to see the names.")

# 2d. Insert the final paragraph describing the bar-graph screenshot.
$xml2 = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r><w:t>There is an image of this bar graph under drawings folder called ' + $q1 + '</w:t></w:r>' +
              '<w:proofErr w:type="spellStart"/>' +
              '<w:r><w:t>pagerank_graph</w:t></w:r>' +
              '<w:proofErr w:type="spellEnd"/>' +
              '<w:r><w:t>' + $q2 + '</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$endOfBookmarkPara = $d.Range($bookmarkPara.Range.End - 1, $bookmarkPara.Range.End - 1)
$endOfBookmarkPara.InsertXML($xml2)

Write-Host ("Final paragraph count=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host ($i.ToString() + ": [" + $d.Paragraphs.Item($i).Range.Text + "]")
}
